# "Add 4K DA for Camera Ret" - populate the Device column (B) for the
# second block of rows (rows 2-26) with the device names that correspond
# to each IP Scope row, widen column B to fit the new text, and leave the
# selection where the author finished editing (H22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in column B with the device names. The order below reproduces the
# order in which the distinct strings were first entered (and therefore
# the order they land in the shared-string table), even though the rows
# themselves are filled out of numeric order.
$ws.Range("B2").Value  = "CCU PC"
$ws.Range("B16").Value = "CCU Mon 02/03"
$ws.Range("B15").Value = "CCU Mon 01"
$ws.Range("B17").Value = "CCU Mon 04/05"
$ws.Range("B18").Value = "CCU Mon 06/07"
$ws.Range("B12").Value = "RTR Panel 01 (CCU)"
$ws.Range("B4").Value  = "ACR PC"
$ws.Range("B14").Value = "Audio Mixer 01"
$ws.Range("B19").Value = "ACR Mon 01"
$ws.Range("B20").Value = "ACR Mon 02/03"
$ws.Range("B13").Value = "RTR Panel 02 (ACR)"
$ws.Range("B23").Value = "VTR 01"
$ws.Range("B24").Value = "VTR 02"
$ws.Range("B21").Value = "VT Mon 01/02"
$ws.Range("B22").Value = "VT Mon 03/04"
$ws.Range("B6").Value  = "CG 02"
$ws.Range("B3").Value  = "VMX PC"
$ws.Range("B10").Value = "VMX 01"
$ws.Range("B8").Value  = "RTR 01"
$ws.Range("B9").Value  = "Frame 01"
$ws.Range("B5").Value  = "CG 01"
$ws.Range("B7").Value  = "Prompt 01"
$ws.Range("B25").Value = "VTR 03"
$ws.Range("B26").Value = "VTR 04"
$ws.Range("B11").Value = "VMX Tub"

# Column B now holds longer device labels such as "RTR Panel 01 (CCU)";
# widen it so the text fits (splits the former merged 2:3 column-width
# entry into its own width for column B).
$ws.Columns("B").ColumnWidth = 16.619791666666668

# Leave the selection where the author ended up working.
$ws.Range("H22").Select() | Out-Null
